# Update the severity-level header labels across every worksheet.
# Old -> New:
#   % 1-2 -> % severity levels 1-2
#   # 1-2 -> # severity levels 1-2
#   % 3   -> % severity level 3
#   # 3   -> # severity level 3
#   % 4   -> % severity level 4
#   # 4   -> # severity level 4
#   % 5   -> % severity level 5
#   # 5   -> # severity level 5

$wb = $excel.ActiveWorkbook

$headerMap = @{
    "% 1-2" = "% severity levels 1-2"
    "# 1-2" = "# severity levels 1-2"
    "% 3"   = "% severity level 3"
    "# 3"   = "# severity level 3"
    "% 4"   = "% severity level 4"
    "# 4"   = "# severity level 4"
    "% 5"   = "% severity level 5"
    "# 5"   = "# severity level 5"
}

foreach ($ws in $wb.Worksheets) {
    # Header labels live in row 1, columns E through L on every sheet.
    for ($col = 5; $col -le 12; $col++) {
        $cell = $ws.Cells.Item(1, $col)
        $current = $cell.Value2
        if ($headerMap.ContainsKey($current)) {
            $cell.Value2 = $headerMap[$current]
        }
    }
}
